$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C1").Value = "ST"
$ws.Range("C10").Value = 152
$ws.Range("C12").Value = 152

$ws.Range("F9").Select()
